$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1499.625
$ws.Range("I19").Value = 1299.6
$ws.Range("J19").Value = 1833
$ws.Range("K19").Value = 1299.6
$ws.Range("L19").Value = 1833
$ws.Range("M19").Value = -1124.6
$ws.Range("N19").Value = -2183

$ws.Range("H62").Value = 7332.6665
$ws.Range("I62").Value = 1999
$ws.Range("J62").Value = 9999.5
$ws.Range("K62").Value = 1999
$ws.Range("L62").Value = 9999.5
$ws.Range("M62").Value = -1375
$ws.Range("N62").Value = -11247.5

$ws.Range("H65").Value = 7332.6665
$ws.Range("I65").Value = 1999
$ws.Range("J65").Value = 9999.5
$ws.Range("K65").Value = 9995
$ws.Range("L65").Value = 49997.5
$ws.Range("M65").Value = -6875
$ws.Range("N65").Value = -56237.5

$ws.Range("H69").Value = 500
$ws.Range("J69").Value = 500
$ws.Range("L69").Value = 1500
$ws.Range("N69").Value = -3248

$ws.Range("H72").Value = 500
$ws.Range("J72").Value = 500
$ws.Range("L72").Value = 4500
$ws.Range("N72").Value = -13236

$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H132").Value = 8091.7144
$ws.Range("I132").Value = 8215.691999999999
$ws.Range("K132").Value = 24647.076
$ws.Range("M132").Value = -22117.076

$ws.Range("H137").Value = 1521.8667
$ws.Range("I137").Value = 1412.5
$ws.Range("J137").Value = 1740.6
$ws.Range("K137").Value = 4237.5
$ws.Range("L137").Value = 5221.799999999999
$ws.Range("M137").Value = -1687.5
$ws.Range("N137").Value = -10321.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 417.16666
$ws.Range("J4").Value = 417.66666
$ws.Range("L4").Value = 417.66666
$ws.Range("N4").Value = -649.66666

$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H110").Value = 669.2
$ws.Range("I110").Value = 685.5714
$ws.Range("J110").Value = 631
$ws.Range("K110").Value = 685.5714
$ws.Range("L110").Value = 631
$ws.Range("M110").Value = 1359.4286
$ws.Range("N110").Value = -4721

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1793.8
$ws.Range("I94").Value = 1793.8
$ws.Range("K94").Value = 1793.8
$ws.Range("M94").Value = -1342.8

$ws.Range("H105").Value = 1390.7778
$ws.Range("I105").Value = 1256.4286
$ws.Range("J105").Value = 1861
$ws.Range("K105").Value = 1256.4286
$ws.Range("L105").Value = 1861
$ws.Range("M105").Value = 490.5714
$ws.Range("N105").Value = -5355

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 24.357143
$ws.Range("I7").Value = 21.11111
$ws.Range("K7").Value = 21.11111
$ws.Range("M7").Value = 91.88889

$ws.Range("H36").Value = 4750
$ws.Range("I36").Value = 4750
$ws.Range("K36").Value = 4750
$ws.Range("M36").Value = -4362

$ws.Range("H40").Value = 4750
$ws.Range("I40").Value = 4750
$ws.Range("K40").Value = 4750
$ws.Range("M40").Value = -4590

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 7.3333335
$ws.Range("I5").Value = 7.3333335
$ws.Range("K5").Value = 22.0000005
$ws.Range("M5").Value = 89.9999995

$ws.Range("H38").Value = 1327.5
$ws.Range("J38").Value = 93.2
$ws.Range("L38").Value = 279.6
$ws.Range("N38").Value = -973.6

$ws.Range("H87").Value = 0
$ws.Range("I87").Value = 0
$ws.Range("K87").Value = 0
$ws.Range("M87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("I90").Value = 0
$ws.Range("K90").Value = 0
$ws.Range("M90").ClearContents()

$ws.Range("H131").Value = 1561.2858
$ws.Range("I131").Value = 1410.6666
$ws.Range("J131").Value = 1674.25
$ws.Range("K131").Value = 4231.9998
$ws.Range("L131").Value = 5022.75
$ws.Range("M131").Value = 808.0002000000004
$ws.Range("N131").Value = -15102.75

$ws.Range("H135").Value = 7.3333335
$ws.Range("I135").Value = 7.3333335
$ws.Range("K135").Value = 66.0000015
$ws.Range("M135").Value = 2468.9999985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 140.83333
$ws.Range("J2").Value = 204.75
$ws.Range("L2").Value = 204.75
$ws.Range("N2").Value = -430.75

$ws.Range("H80").Value = 15673.5
$ws.Range("J80").Value = 49900
$ws.Range("L80").Value = 49900
$ws.Range("N80").Value = -51896

$ws.Range("H83").Value = 15673.5
$ws.Range("J83").Value = 49900
$ws.Range("L83").Value = 249500
$ws.Range("N83").Value = -259484

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5566.1665
$ws.Range("I16").Value = 5566.1665
$ws.Range("K16").Value = 5566.1665
$ws.Range("M16").Value = -5396.1665

$ws.Range("H22").Value = 887.375
$ws.Range("I22").Value = 575
$ws.Range("J22").Value = 1199.75
$ws.Range("K22").Value = 575
$ws.Range("L22").Value = 1199.75
$ws.Range("M22").Value = -280
$ws.Range("N22").Value = -1789.75

$ws.Range("H27").Value = 887.375
$ws.Range("I27").Value = 575
$ws.Range("J27").Value = 1199.75
$ws.Range("K27").Value = 575
$ws.Range("L27").Value = 1199.75
$ws.Range("M27").Value = -468
$ws.Range("N27").Value = -1413.75

$ws.Range("H30").Value = 1126
$ws.Range("I30").Value = 1315
$ws.Range("J30").Value = 1000
$ws.Range("K30").Value = 1315
$ws.Range("L30").Value = 1000
$ws.Range("M30").Value = -1207
$ws.Range("N30").Value = -1216

$ws.Range("H40").Value = 3900
$ws.Range("I40").Value = 3900
$ws.Range("K40").Value = 3900
$ws.Range("M40").Value = -3764

$ws.Range("H55").Value = 1674
$ws.Range("I55").Value = 608
$ws.Range("J55").Value = 2740
$ws.Range("K55").Value = 608
$ws.Range("L55").Value = 2740
$ws.Range("M55").Value = -435
$ws.Range("N55").Value = -3086

$ws.Range("H93").Value = 2000
$ws.Range("I93").Value = 2000
$ws.Range("K93").Value = 2000
$ws.Range("M93").Value = -752

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 108333.336
$ws.Range("J123").Value = 108333.336
$ws.Range("L123").Value = 108333.336
$ws.Range("N123").Value = -118133.336

$ws.Range("H132").Value = 960
$ws.Range("I132").Value = 960
$ws.Range("K132").Value = 2880
$ws.Range("M132").Value = -350
